# Update the price-list date shown in cell A1 (Hoja1) from 45344 (2024-02-22)
# to 45405 (2024-04-23), as part of adding a lock screen in step 1 to
# mitigate multiple requests.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45405

# Re-create the merged ranges so they are re-written in the same order
# Excel produces after the edit.
$ws.Range("A1:D1").UnMerge()
$ws.Range("A12:D12").UnMerge()
$ws.Range("B19:C19").UnMerge()
$ws.Range("A13:D13").UnMerge()
$ws.Range("A11:D11").UnMerge()
$ws.Range("B18:C18").UnMerge()
$ws.Range("B17:C17").UnMerge()
$ws.Range("A14:D14").UnMerge()

$ws.Range("A1:D1").Merge()
$ws.Range("A12:D12").Merge()
$ws.Range("B19:C19").Merge()
$ws.Range("A13:D13").Merge()
$ws.Range("A11:D11").Merge()
$ws.Range("B18:C18").Merge()
$ws.Range("B17:C17").Merge()
$ws.Range("A14:D14").Merge()
